$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.921.40"
$ws.Range("E2").Value = "  -1.02%  "

$ws.Range("D3").Value = "1.833.31"
$ws.Range("E3").Value = "  -1.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.30"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6918"
$ws.Range("E6").Value = "  -1.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07683"
$ws.Range("E8").Value = "  -1.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3045"
$ws.Range("E9").Value = "  -2.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.32"
$ws.Range("E10").Value = "  -3.03%  "

$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "93.32"
$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("D13").Value = "1.829.22"
$ws.Range("E13").Value = "  -1.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.101"
$ws.Range("E14").Value = "  -0.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6816"
$ws.Range("E15").Value = "  -1.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.550"
$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008253"
$ws.Range("E17").Value = "  -3.48%  "

$ws.Range("D18").Value = "28.930.63"
$ws.Range("E18").Value = "  -1.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.60"
$ws.Range("E19").Value = "  -3.12%  "

$ws.Range("D20").Value = "2.075.20"
$ws.Range("E20").Value = "  -1.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.69"
$ws.Range("E21").Value = "  -1.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.458"
$ws.Range("E23").Value = "  -1.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9991"
$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1498"
$ws.Range("E25").Value = "  -2.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.95"
$ws.Range("E26").Value = "  -1.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.738"
$ws.Range("E27").Value = "  -2.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.17"
$ws.Range("E28").Value = "  -2.34%  "

$ws.Range("E29").Value = "  -2.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.222"
$ws.Range("E30").Value = "  -1.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.136"
$ws.Range("E31").Value = "  -2.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.193"
$ws.Range("E32").Value = "  -1.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05107"
$ws.Range("E33").Value = "  -2.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7770"
$ws.Range("E34").Value = "  +2.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.854"
$ws.Range("E35").Value = "  -1.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.144"
$ws.Range("E36").Value = "  -2.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.694"
$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("D38").Value = "1.284.69"
$ws.Range("E38").Value = "  +3.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01859"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.698"
$ws.Range("E40").Value = "  -1.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9502"
$ws.Range("E41").Value = "  +5.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.159"
$ws.Range("E42").Value = "  +4.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.77"
$ws.Range("E43").Value = "  -2.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9989"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.678"
$ws.Range("E45").Value = "  +1.52%  "

$ws.Range("E46").Value = "  -1.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5162"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("D48").Value = "1.974.54"
$ws.Range("E48").Value = "  -1.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.78"
$ws.Range("E49").Value = "  -6.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.753"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.953"
$ws.Range("E51").Value = "  -0.98%  "
